$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "even_MAG-GUT89162.fa"
$ws.Range("A3").Value = "even_MAG-GUT89246.fa"
$ws.Range("A4").Value = "even_MAG-GUT89784.fa"
$ws.Range("A5").Value = "even_MAG-GUT89815.fa"
$ws.Range("A6").Value = "even_MAG-GUT89852.fa"
$ws.Range("A7").Value = "even_MAG-GUT90020.fa"
$ws.Range("A8").Value = "even_MAG-GUT90054.fa"
$ws.Range("A9").Value = "even_MAG-GUT90362.fa"
$ws.Range("A10").Value = "even_MAG-GUT90441.fa"
$ws.Range("A11").Value = "even_MAG-GUT90682.fa"
$ws.Range("A12").Value = "even_MAG-GUT90963.fa"
$ws.Range("A13").Value = "even_MAG-GUT91014.fa"

$ws.Range("A14:D18").EntireRow.Delete()
